$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 224
$ws.Range("I2").Value = 208.2
$ws.Range("K2").Value = 208.2
$ws.Range("M2").Value = -95.19999999999999
$ws.Range("H9").Value = 122.63636
$ws.Range("I9").Value = 56.25
$ws.Range("K9").Value = 56.25
$ws.Range("M9").Value = 112.75
$ws.Range("H38").Value = 746
$ws.Range("H58").Value = 2724.75
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 2724.75
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 8174.25
$ws.Range("M58").Value = $null
$ws.Range("N58").Value = -8474.25
$ws.Range("H116").Value = 3370.6553
$ws.Range("I116").Value = 2533.4707
$ws.Range("J116").Value = 4556.6665
$ws.Range("K116").Value = 2533.4707
$ws.Range("L116").Value = 4556.6665
$ws.Range("M116").Value = 908.5293000000001
$ws.Range("N116").Value = -11440.6665
$ws.Range("H132").Value = 14650.194
$ws.Range("I132").Value = 16926.475
$ws.Range("J132").Value = 2027.1818
$ws.Range("K132").Value = 50779.425
$ws.Range("L132").Value = 6081.5454
$ws.Range("M132").Value = -48249.425
$ws.Range("N132").Value = -11141.5454
$ws.Range("H137").Value = 1240.9375
$ws.Range("I137").Value = 1256.5883
$ws.Range("J137").Value = 1223.2
$ws.Range("K137").Value = 3769.7649
$ws.Range("L137").Value = 3669.6
$ws.Range("M137").Value = -1219.7649
$ws.Range("N137").Value = -8769.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 29000
$ws.Range("J7").Value = 29000
$ws.Range("L7").Value = 29000
$ws.Range("N7").Value = -29228
$ws.Range("H32").Value = 12452.359
$ws.Range("I32").Value = 10179.135
$ws.Range("J32").Value = 54507
$ws.Range("K32").Value = 10179.135
$ws.Range("L32").Value = 54507
$ws.Range("M32").Value = -9892.135
$ws.Range("N32").Value = -55081
$ws.Range("H52").Value = 40597.5
$ws.Range("J52").Value = 40597.5
$ws.Range("L52").Value = 40597.5
$ws.Range("N52").Value = -41233.5
$ws.Range("H61").Value = 2397.9092
$ws.Range("I61").Value = 2215.875
$ws.Range("K61").Value = 2215.875
$ws.Range("M61").Value = -2003.875
$ws.Range("H74").Value = 1680.6136
$ws.Range("I74").Value = 1089.2059
$ws.Range("J74").Value = 3691.4
$ws.Range("K74").Value = 1089.2059
$ws.Range("L74").Value = 3691.4
$ws.Range("M74").Value = -215.2058999999999
$ws.Range("N74").Value = -5439.4
$ws.Range("H77").Value = 1680.6136
$ws.Range("I77").Value = 1089.2059
$ws.Range("J77").Value = 3691.4
$ws.Range("K77").Value = 5446.0295
$ws.Range("L77").Value = 18457
$ws.Range("M77").Value = -1078.0295
$ws.Range("N77").Value = -27193
$ws.Range("H136").Value = 2397.9092
$ws.Range("I136").Value = 2215.875
$ws.Range("K136").Value = 6647.625
$ws.Range("M136").Value = -4097.625

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1749.119
$ws.Range("I134").Value = 1350.5667
$ws.Range("J134").Value = 2745.5
$ws.Range("K134").Value = 4051.7001
$ws.Range("L134").Value = 8236.5
$ws.Range("M134").Value = -1516.7001
$ws.Range("N134").Value = -13306.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 12925300
$ws.Range("I6").Value = 25850000
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 25850000
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -25849887
$ws.Range("N6").Value = -826
$ws.Range("H7").Value = 639.9474
$ws.Range("I7").Value = 777.93335
$ws.Range("J7").Value = 122.5
$ws.Range("K7").Value = 777.93335
$ws.Range("L7").Value = 122.5
$ws.Range("M7").Value = -664.93335
$ws.Range("N7").Value = -348.5
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("M17").Value = -826
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("N25").Value = $null
$ws.Range("H58").Value = 2100.2104
$ws.Range("I58").Value = 1262.5834
$ws.Range("J58").Value = 3536.1428
$ws.Range("K58").Value = 1262.5834
$ws.Range("L58").Value = 3536.1428
$ws.Range("M58").Value = -1059.5834
$ws.Range("N58").Value = -3942.1428
$ws.Range("H59").Value = 31266.666
$ws.Range("I59").Value = 20000
$ws.Range("J59").Value = 33520
$ws.Range("K59").Value = 20000
$ws.Range("L59").Value = 33520
$ws.Range("M59").Value = -18855
$ws.Range("N59").Value = -35810
$ws.Range("H132").Value = 1002261.06
$ws.Range("I132").Value = 1667726.5
$ws.Range("J132").Value = 4062.9
$ws.Range("K132").Value = 5003179.5
$ws.Range("L132").Value = 12188.7
$ws.Range("M132").Value = -5000649.5
$ws.Range("N132").Value = -17248.7
$ws.Range("H136").Value = 2100.2104
$ws.Range("I136").Value = 1262.5834
$ws.Range("J136").Value = 3536.1428
$ws.Range("K136").Value = 3787.7502
$ws.Range("L136").Value = 10608.4284
$ws.Range("M136").Value = -1237.7502
$ws.Range("N136").Value = -15708.4284
$ws.Range("H137").Value = 28200
$ws.Range("J137").Value = 28200
$ws.Range("L137").Value = 28200
$ws.Range("N137").Value = -38400

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 778.21
$ws.Range("I131").Value = 461.88235
$ws.Range("K131").Value = 1385.64705
$ws.Range("M131").Value = 3654.35295

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6269.4287
$ws.Range("I70").Value = 4810.769
$ws.Range("J70").Value = 8639.75
$ws.Range("K70").Value = 4810.769
$ws.Range("L70").Value = 8639.75
$ws.Range("M70").Value = -4540.769
$ws.Range("N70").Value = -9179.75
$ws.Range("H73").Value = 6269.4287
$ws.Range("I73").Value = 4810.769
$ws.Range("J73").Value = 8639.75
$ws.Range("K73").Value = 4810.769
$ws.Range("L73").Value = 8639.75
$ws.Range("M73").Value = -3874.769
$ws.Range("N73").Value = -10511.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = $null
$ws.Range("N56").Value = $null
$ws.Range("H136").Value = 4588.467
$ws.Range("I136").Value = 5195.1
$ws.Range("J136").Value = 3375.2
$ws.Range("K136").Value = 15585.3
$ws.Range("L136").Value = 10125.6
$ws.Range("M136").Value = -13035.3
$ws.Range("N136").Value = -15225.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5981.294
$ws.Range("I136").Value = 4186.9165
$ws.Range("J136").Value = 7576.2964
$ws.Range("K136").Value = 12560.7495
$ws.Range("L136").Value = 22728.8892
$ws.Range("M136").Value = -10010.7495
$ws.Range("N136").Value = -27828.8892
